$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.639.80'
$ws.Range("E2").Value = '  -1.32%  '
Set-TextValue $ws.Range("D3") '2.449.87'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '556.06'
$ws.Range("E5").Value = '  -2.40%  '
Set-TextValue $ws.Range("D6") '160.88'
$ws.Range("E6").Value = '  -2.61%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.92%  '
Set-TextValue $ws.Range("D9") '2.448.74'
$ws.Range("E9").Value = '  -1.73%  '
Set-TextValue $ws.Range("D10") '0.148'
$ws.Range("E10").Value = '  -7.30%  '
$ws.Range("E11").Value = '  -0.94%  '
Set-TextValue $ws.Range("D12") '0.333'
$ws.Range("E12").Value = '  -5.29%  '
Set-TextValue $ws.Range("D13") '4.78'
$ws.Range("E13").Value = '  -2.03%  '
Set-TextValue $ws.Range("D14") '2.891.37'
$ws.Range("E14").Value = '  -1.87%  '
Set-TextValue $ws.Range("D15") '68.509.68'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("E16").Value = '  -4.59%  '
Set-TextValue $ws.Range("D17") '23.27'
$ws.Range("E17").Value = '  -3.79%  '
Set-TextValue $ws.Range("D18") '2.450.99'
$ws.Range("E18").Value = '  -1.77%  '
Set-TextValue $ws.Range("D19") '10.70'
$ws.Range("E19").Value = '  -4.23%  '
Set-TextValue $ws.Range("D20") '339.93'
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("E21").Value = '  -5.82%  '
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("E24").Value = '  +0.30%  '
Set-TextValue $ws.Range("D25") '1.84'
$ws.Range("E25").Value = '  -2.69%  '
Set-TextValue $ws.Range("D26") '66.38'
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("E27").Value = '  -6.01%  '
$ws.Range("E28").Value = '  -2.25%  '
Set-TextValue $ws.Range("D29") '1.00'
$ws.Range("E29").Value = '  +0.06%  '
Set-TextValue $ws.Range("D30") '8.06'
$ws.Range("E30").Value = '  -6.14%  '
Set-TextValue $ws.Range("D31") '0.0₃0810'
$ws.Range("E31").Value = '  -6.98%  '
$ws.Range("E32").Value = '  -6.24%  '
Set-TextValue $ws.Range("D33") '434.16'
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("E35").Value = '  -5.65%  '
$ws.Range("E36").Value = '  -5.91%  '
Set-TextValue $ws.Range("D37") '156.17'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  -3.33%  '
Set-TextValue $ws.Range("D41") '17.79'
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  -4.25%  '
Set-TextValue $ws.Range("D44") '37.40'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("E45").Value = '  -7.82%  '
$ws.Range("E46").Value = '  +1.77%  '
Set-TextValue $ws.Range("D47") '2.03'
$ws.Range("E47").Value = '  -6.17%  '
Set-TextValue $ws.Range("D48") '131.77'
$ws.Range("E48").Value = '  -4.97%  '
Set-TextValue $ws.Range("D49") '3.32'
$ws.Range("E49").Value = '  -2.95%  '
Set-TextValue $ws.Range("D50") '0.0713'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("E51").Value = '  -4.91%  '
